$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.878.02"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "3.378.30"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.70"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.26%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "187.30"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.49%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +2.08%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.185"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  +1.14%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "47.61"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("D13").Value = "3.918.47"
$ws.Range("E13").Value = "  +0.81%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "640.52"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +7.69%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "67.805.20"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "3.378.80"
$ws.Range("E18").Value = "  +1.06%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "18.10"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.17"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.912"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.23%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "18.01"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  +1.71%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "99.79"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +4.92%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.81"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.17%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "32.72"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.23%  "
$ws.Range("E29").Value = "  +2.26%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.94"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.40%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "613.72"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.20%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.87"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").Value = "4.029.76"
$ws.Range("E33").Value = "  +6.93%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "11.14"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("E36").Value = "  +0.00%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "56.33"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("E38").Value = "  +6.38%  "
$ws.Range("E39").Value = "  +3.88%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "33.85"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  +1.30%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0424"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +11.00%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "128.04"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.97%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.77"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.91%  "
